# experienceCleanUp: Minor housekeeping on some formatting and adding some
# things for additional experience

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. WorkExperience sheet: insert an "is_additional" column between
#    is_current and description, add a 4th (additional) job row, and
#    a handful of pre-formatted blank rows below it.
# ---------------------------------------------------------------------
$wsWork = $wb.Worksheets.Item("WorkExperience")

# Shift the existing "description" (F) / "accomplishments" (G) columns
# one slot to the right -> G / H, working from the rightmost column
# back so we never clobber data before it has been read. The header
# row also carries its bold/banner formatting (style s="3") along, so
# copy formats there first.
$wsWork.Cells.Item(1, 7).Copy()
$wsWork.Cells.Item(1, 8).PasteSpecial(-4122)
$wsWork.Cells.Item(1, 6).Copy()
$wsWork.Cells.Item(1, 7).PasteSpecial(-4122)

for ($r = 1; $r -le 3; $r++) {
    $gVal = $wsWork.Cells.Item($r, 7).Value()
    $fVal = $wsWork.Cells.Item($r, 6).Value()
    $wsWork.Cells.Item($r, 8).Value = $gVal
    $wsWork.Cells.Item($r, 7).Value = $fVal
}

# F1 becomes the new "is_additional" header (still styled s="3" from
# the copy/paste above).
$wsWork.Cells.Item(1, 6).Value = "is_additional"

# New boolean "is_additional" data cells for the two existing rows
# (both FALSE), styled with an explicit (if redundant) font so that a
# dedicated cell style gets produced.
$wsWork.Cells.Item(2, 6).Value = $false
$wsWork.Cells.Item(2, 6).Font.Name = "Calibri"

$wsWork.Cells.Item(3, 6).Value = $false
$wsWork.Cells.Item(3, 6).Font.Name = "Calibri"

# Narrow column F down from the old "description" width to fit the new
# boolean values/header.
$wsWork.Columns.Item(6).ColumnWidth = 10.17

# Row 4: a new "additional experience" entry.
$wsWork.Cells.Item(4, 1).Value = "My First job"
$wsWork.Cells.Item(4, 2).Value = "Doesn’t matter"
$wsWork.Cells.Item(4, 3).Value = "2016-01"
$wsWork.Cells.Item(4, 4).Value = "2018-05"
$wsWork.Cells.Item(4, 5).Value = $false
$wsWork.Cells.Item(4, 6).Value = $true
$wsWork.Cells.Item(4, 6).Font.Name = "Calibri"

# Rows 5-8: pre-formatted (but otherwise empty) continuation cells in
# column F, matching the style used for the is_additional checkboxes.
$wsWork.Cells.Item(2, 6).Copy()
$wsWork.Cells.Item(5, 6).PasteSpecial(-4122)
$wsWork.Cells.Item(6, 6).PasteSpecial(-4122)
$wsWork.Cells.Item(7, 6).PasteSpecial(-4122)
$wsWork.Cells.Item(8, 6).PasteSpecial(-4122)
$wsWork.Cells.Item(5, 6).ClearContents()
$wsWork.Cells.Item(6, 6).ClearContents()
$wsWork.Cells.Item(7, 6).ClearContents()
$wsWork.Cells.Item(8, 6).ClearContents()

# Selection / active-cell bookkeeping to match the saved view state.
$wsWork.Range("E16").Select()

# ---------------------------------------------------------------------
# 2. Sheet-tab bookkeeping: WorkExperience becomes the active/selected
#    tab instead of Education.
# ---------------------------------------------------------------------
$wsWork.Activate()
